# "Moved app.py and all files to root" - append the latest login-history
# row (row 15) to the user log sheet, matching the app's new log append.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A15").Value = "kumarshashwat890@gmail.com"
$ws.Range("B15").Value = "Shashwat kumar"
$ws.Range("C15").Value = "2025-07-12 17:10:01"

# Logout Time (D15) has not happened yet, so the source log writes an
# explicit empty string rather than leaving the cell truly blank. Assigning
# "" directly clears/removes the cell instead of storing an empty string,
# so write a text value and strip the apostrophe text-prefix style it adds.
$ws.Range("D15").Value = "'"
$ws.Range("D15").Style = "Normal"
